# Auto-generated Excel COM-interop script to apply BRVM recommendation refresh
$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# --- Update "Recommandations" sheet (rows 2-48 keep data; row 49 removed) ---
# Row 2
$wsReco.Cells.Item(2, 1).Value = "SAFCA CI"
$wsReco.Cells.Item(2, 3).Value = 4
$wsReco.Cells.Item(2, 4).Value = 3325
$wsReco.Cells.Item(2, 5).Value = 920
# Row 3
$wsReco.Cells.Item(3, 1).Value = "BRVM - SERVICES PUBLICS"
$wsReco.Cells.Item(3, 3).Value = 8
$wsReco.Cells.Item(3, 4).Value = 3245.54
$wsReco.Cells.Item(3, 5).Value = 101.52
# Row 4
$wsReco.Cells.Item(4, 1).Value = "CFAO MOTORS CI"
$wsReco.Cells.Item(4, 4).Value = 2625
$wsReco.Cells.Item(4, 5).Value = 650
# Row 5
$wsReco.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$wsReco.Cells.Item(5, 4).Value = 2614.33
$wsReco.Cells.Item(5, 5).Value = 656
# Row 6
$wsReco.Cells.Item(6, 1).Value = "NEI-CEDA CI"
$wsReco.Cells.Item(6, 4).Value = 2340
# Row 7
$wsReco.Cells.Item(7, 1).Value = "SETAO CI"
$wsReco.Cells.Item(7, 4).Value = 2340
$wsReco.Cells.Item(7, 5).Value = 590
# Row 8
$wsReco.Cells.Item(8, 4).Value = 2275
$wsReco.Cells.Item(8, 5).Value = 570
# Row 9
$wsReco.Cells.Item(9, 4).Value = 2075
$wsReco.Cells.Item(9, 5).Value = 520
# Row 10
$wsReco.Cells.Item(10, 4).Value = 1443.32
$wsReco.Cells.Item(10, 5).Value = 361.94
# Row 11
$wsReco.Cells.Item(11, 4).Value = 1392.75
$wsReco.Cells.Item(11, 5).Value = 345.13
# Row 12
$wsReco.Cells.Item(12, 4).Value = 1261.38
$wsReco.Cells.Item(12, 5).Value = 314.38
# Row 13
$wsReco.Cells.Item(13, 4).Value = 1028.61
$wsReco.Cells.Item(13, 5).Value = 261.56
# Row 14
$wsReco.Cells.Item(14, 4).Value = 848.76
$wsReco.Cells.Item(14, 5).Value = 215.03
# Row 15
$wsReco.Cells.Item(15, 4).Value = 750.38
$wsReco.Cells.Item(15, 5).Value = 187.73
# Row 16
$wsReco.Cells.Item(16, 4).Value = 550.75
$wsReco.Cells.Item(16, 5).Value = 137.14
# Row 17
$wsReco.Cells.Item(17, 4).Value = 515.99
$wsReco.Cells.Item(17, 5).Value = 128.2
# Row 18
$wsReco.Cells.Item(18, 4).Value = 486.95
$wsReco.Cells.Item(18, 5).Value = 120.6
# Row 19
$wsReco.Cells.Item(19, 4).Value = 478.57
$wsReco.Cells.Item(19, 5).Value = 118.52
# Row 20
$wsReco.Cells.Item(20, 1).Value = "BRVM - ENERGIE"
$wsReco.Cells.Item(20, 4).Value = 424.56
$wsReco.Cells.Item(20, 5).Value = 107.37
# Row 21
$wsReco.Cells.Item(21, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsReco.Cells.Item(21, 4).Value = 423.44
$wsReco.Cells.Item(21, 5).Value = 105.17
# Row 22
$wsReco.Cells.Item(22, 4).Value = 369.11
$wsReco.Cells.Item(22, 5).Value = 92.36
# Row 23
$wsReco.Cells.Item(23, 1).Value = "SAFCA CI (SAFC)"
$wsReco.Cells.Item(23, 4).Value = 29.05
$wsReco.Cells.Item(23, 5).Value = 6.98
# Row 24
$wsReco.Cells.Item(24, 1).Value = "UNILEVER CI (UNLC)"
$wsReco.Cells.Item(24, 4).Value = 22.38
$wsReco.Cells.Item(24, 5).Value = 7.49
# Row 25
$wsReco.Cells.Item(25, 1).Value = "UNIWAX CI (UNXC)"
$wsReco.Cells.Item(25, 2).Value = 2
$wsReco.Cells.Item(25, 4).Value = 5.35
$wsReco.Cells.Item(25, 5).Value = 1.75
# Row 26
# Row 27
# Row 28
# Row 29
# Row 30
$wsReco.Cells.Item(30, 1).Value = "BANK OF AFRICA BN (BOAB)"
$wsReco.Cells.Item(30, 2).Value = 1
$wsReco.Cells.Item(30, 3).Value = 0
$wsReco.Cells.Item(30, 4).Value = 2.86
$wsReco.Cells.Item(30, 5).Value = 2.86
$wsReco.Cells.Item(30, 7).Value = "➖ Neutre"
# Row 31
$wsReco.Cells.Item(31, 1).Value = "SICOR CI (SICC)"
$wsReco.Cells.Item(31, 4).Value = 2.8
$wsReco.Cells.Item(31, 5).Value = 2.8
# Row 32
$wsReco.Cells.Item(32, 1).Value = "FILTISAC CI (FTSC)"
$wsReco.Cells.Item(32, 4).Value = 2.09
$wsReco.Cells.Item(32, 5).Value = -0.8
# Row 33
$wsReco.Cells.Item(33, 1).Value = "NEI-CEDA CI (NEIC)"
$wsReco.Cells.Item(33, 4).Value = 0.85
$wsReco.Cells.Item(33, 5).Value = 1.69
# Row 34
$wsReco.Cells.Item(34, 1).Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$wsReco.Cells.Item(34, 4).Value = 0.1
$wsReco.Cells.Item(34, 5).Value = -3.03
# Row 35
# Row 36
$wsReco.Cells.Item(36, 1).Value = "BERNABE CI (BNBC)"
$wsReco.Cells.Item(36, 4).Value = -0.26
$wsReco.Cells.Item(36, 5).Value = -3.85
# Row 37
# Row 38
$wsReco.Cells.Item(38, 1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$wsReco.Cells.Item(38, 4).Value = -1.6
$wsReco.Cells.Item(38, 5).Value = -1.6
# Row 39
$wsReco.Cells.Item(39, 1).Value = "BANK OF AFRICA NG (BOAN)"
$wsReco.Cells.Item(39, 2).Value = 1
$wsReco.Cells.Item(39, 3).Value = 2
$wsReco.Cells.Item(39, 4).Value = -1.95
$wsReco.Cells.Item(39, 5).Value = 2.24
$wsReco.Cells.Item(39, 7).Value = "👀 À surveiller"
# Row 40
$wsReco.Cells.Item(40, 1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$wsReco.Cells.Item(40, 4).Value = -2.22
$wsReco.Cells.Item(40, 5).Value = -2.22
# Row 41
$wsReco.Cells.Item(41, 1).Value = "SITAB CI (STBC)"
$wsReco.Cells.Item(41, 4).Value = -2.5
$wsReco.Cells.Item(41, 5).Value = -2.5
# Row 42
$wsReco.Cells.Item(42, 1).Value = "BICI CI (BICC)"
$wsReco.Cells.Item(42, 4).Value = -2.85
$wsReco.Cells.Item(42, 5).Value = -2.85
# Row 43
$wsReco.Cells.Item(43, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$wsReco.Cells.Item(43, 4).Value = -3.42
$wsReco.Cells.Item(43, 5).Value = -3.42
# Row 44
$wsReco.Cells.Item(44, 1).Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Cells.Item(44, 4).Value = -3.64
$wsReco.Cells.Item(44, 5).Value = -3.64
# Row 45
$wsReco.Cells.Item(45, 1).Value = "CFAO MOTORS CI (CFAC)"
$wsReco.Cells.Item(45, 4).Value = -3.65
$wsReco.Cells.Item(45, 5).Value = -3.65
# Row 46
$wsReco.Cells.Item(46, 1).Value = "SUCRIVOIRE (SCRC)"
$wsReco.Cells.Item(46, 2).Value = 0
$wsReco.Cells.Item(46, 4).Value = -3.85
$wsReco.Cells.Item(46, 5).Value = -3.85
$wsReco.Cells.Item(46, 7).Value = "➖ Neutre"
# Row 47
$wsReco.Cells.Item(47, 1).Value = "ONATEL BF (ONTBF)"
$wsReco.Cells.Item(47, 3).Value = 2
$wsReco.Cells.Item(47, 4).Value = -4.21
$wsReco.Cells.Item(47, 5).Value = -2.13
# Row 48
$wsReco.Cells.Item(48, 1).Value = "CIE CI (CIEC)"
$wsReco.Cells.Item(48, 3).Value = 2
$wsReco.Cells.Item(48, 4).Value = -8.97
$wsReco.Cells.Item(48, 5).Value = -3.91

# Remove obsolete row 49 (BANK OF AFRICA NG (BOAN) duplicate tail entry)
$wsReco.Range("A49").EntireRow.Delete()

# --- Update "Top_YTD" sheet (rows 2-11) ---
# Row 2
$wsYtd.Cells.Item(2, 2).Value = 7315290.28
# Row 3
$wsYtd.Cells.Item(3, 2).Value = 744581.6
# Row 4
$wsYtd.Cells.Item(4, 1).Value = "CFAO MOTORS CI"
$wsYtd.Cells.Item(4, 2).Value = 326909.38
# Row 5
$wsYtd.Cells.Item(5, 1).Value = "BRVM - AUTRES SECTEURS"
$wsYtd.Cells.Item(5, 2).Value = 322355.92
# Row 6
$wsYtd.Cells.Item(6, 2).Value = 220048.64
# Row 7
$wsYtd.Cells.Item(7, 2).Value = 219941.52
# Row 8
$wsYtd.Cells.Item(8, 2).Value = 199851.03
# Row 9
$wsYtd.Cells.Item(9, 2).Value = 146423.44
# Row 10
$wsYtd.Cells.Item(10, 2).Value = 44997.36
# Row 11
$wsYtd.Cells.Item(11, 2).Value = 40241.9
